$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232) - entire row shift up
$ws.Rows.Item(26).Delete()

# After the above deletion, what was row 28 (SC 92) is now at row 27
$ws.Rows.Item(27).Delete()

# Now fix the B-column "missing data" swap:
# New row 29 (SC 119) loses its B value (becomes missing) - keep an empty string
# to match how other "missing" B cells in the sheet are represented
$ws.Cells.Item(29, 2).Value = ""

# New row 33 (SC 232) gains a B value of -19.5
$ws.Cells.Item(33, 2).Value = -19.5
